# Journal de bord: log the "end of first version" milestone.
# Adds one new row to the table (D15:E15) - a date and an event label -
# reusing the same date-format style already used by the row above it,
# and appending the new label text to the shared-string table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-cell formatting (number format + alignment) from the row
# above so the new date cell picks up the existing style (s="3") instead
# of Excel fabricating a brand-new style entry.
$ws.Range("D14").Copy()
$ws.Range("D15").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new log entry.
$ws.Range("D15").Value = 44266
$ws.Range("E15").Value = "Fin de la première version"

# Leave the selection where the author's cursor ended up.
$ws.Range("T13").Select() | Out-Null
